$d = $word.ActiveDocument
$xml = '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:p>
      <w:pPr>
        <w:jc w:val="center"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Final Summary</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="720"/>
      </w:pPr>
      <w:r>
        <w:t>With a fixed arrival rate of 10 patients per hour the 1 doctor, 1 nurse system was much slower than the 1 doctor, 2 nurse and 2 doctor, 1 nurse systems. Every simulation served an equal number of patients in the week duration of about 1700 people. However, the average wait time for the 1 doctor, 1 nurse system was around 10.58 minutes and the total average visit time was around 18.99 minutes. Where as the average wait time for the 1</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>doctor, 2 nurse system was 3.60 minutes and the average visit time was at 11.46 minutes. I</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> expected the 1 doctor, 2 nurse system to be the most efficient with a fixed number of 10 patients per hour. However, the 2 doctor, 1 nurse system came out to be the fastest with an average wait time of 1.45 minutes but and average visit time of 10.60 minutes. This made me realize that this system was probably the fastest because </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">only </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">the doctors are able to treat patients with an illness level of 11 and above, so there is probably just a slight delay in the other system to treat </w:t>
      </w:r>
      <w:r>
        <w:t>the</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> patient</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">s with higher illness levels </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">that makes it a little longer. The evidence for this is that the wait time for the patients is shorter for the 2 doctor, 1 nurse system but the overall visit time is very close, so the </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">2 doctor, 1 nurse </w:t>
      </w:r>
      <w:r>
        <w:t>emergency room is able to keep up with the patients coming in even though it takes a doctor longer to see a patient. To test this theory, I ran the test again with a fixed number of 20 patients.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> This would test whether my theory would be right because the 2 doctor, 1 nurse system would not be able to keep up with the number of patients coming in because doctors take to long to treat their patients.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> The 1 doctor, 1 nurse system was again not even close and only treated</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> a little over</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> half the </w:t>
      </w:r>
      <w:r>
        <w:t>3400</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> patients that entered the emergency room. The 1 doctor, 2 nurse system was the fastest this time and by a large amount. It had an average wait time of </w:t>
      </w:r>
      <w:r>
        <w:t>19.90 minutes and an average visit time of 27.53 minutes. Whereas the 2 doctor, 1 nurse system had an average wait time of 35.07 minutes and an average visit time of 43.57 minutes. These wait times and visit times are much slower, which supports my theory.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Graph:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Display a plot of the visit time for increasing patient arrival rates, for a combination of doctors and nurses of your choice.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>Fixed arrival rate: 20 patients per hour</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>1 doctor and 1 nurse</w:t>
      </w:r>
      <w:r>
        <w:t>; patients served 2152/3396</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">; </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">1 doctor and 2 </w:t>
      </w:r>
      <w:r>
        <w:t>nurses</w:t>
      </w:r>
      <w:r>
        <w:t>; patient</w:t>
      </w:r>
      <w:r>
        <w:t>s</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> served</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> 3247/3422; average wait time: 19.90 minutes; average visit time: 27.53 minutes.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t>2 doctors and 1 nurse</w:t>
      </w:r>
      <w:r>
        <w:t>; patients served 3314/3179; average wait time: 35.07 minutes; average visit time: 43.57 minutes.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Fixed arrival rate: 10 patients per hour</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">1 doctor and 1 nurse; patients served </w:t>
      </w:r>
      <w:r>
        <w:t>1740</w:t>
      </w:r>
      <w:r>
        <w:t>/</w:t>
      </w:r>
      <w:r>
        <w:t>1744</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">; average wait time: </w:t>
      </w:r>
      <w:r>
        <w:t>10.58</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> minutes; average visit time: </w:t>
      </w:r>
      <w:r>
        <w:t>18</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:t>99</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> minutes.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve">1 doctor and 2 nurses; patients served </w:t>
      </w:r>
      <w:r>
        <w:t>1654</w:t>
      </w:r>
      <w:r>
        <w:t>/</w:t>
      </w:r>
      <w:r>
        <w:t>1657</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">; average wait time: </w:t>
      </w:r>
      <w:r>
        <w:t>3.60</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> minutes; average visit time: </w:t>
      </w:r>
      <w:r>
        <w:t>11.46</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> minutes.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="1"/>
        </w:numPr>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">2 doctors and 1 nurse; patients served </w:t>
      </w:r>
      <w:r>
        <w:t>1651/1652</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">; average wait time: </w:t>
      </w:r>
      <w:r>
        <w:t>1.45</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> minutes; average visit time: </w:t>
      </w:r>
      <w:r>
        <w:t>10.60</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> minutes.</w:t>
      </w:r>
    </w:p></w:document>'
$d.Content.InsertXML($xml)
Write-Output "Paragraph count: $($d.Paragraphs.Count)"
